# Update FRED WTREGEN data: prepend 5 new weekly observations, append 2 new
# weekly observations, and refresh the SeriesInfo metadata to match the
# latest pull from FRED.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# 1) Shift the existing 112 data rows (old rows 2..113) down by 5 rows so
#    there is room for 5 new rows at the top. Walk bottom-up so a row is
#    always copied before its old location is overwritten.
# ---------------------------------------------------------------------------
for ($r = 113; $r -ge 2; $r--) {
    $dstRow = $r + 5
    $src = $ws.Range("A" + $r + ":B" + $r)
    $dst = $ws.Range("A" + $dstRow + ":B" + $dstRow)
    $src.Copy($dst)
}

# ---------------------------------------------------------------------------
# 2) Write the 5 new observations into the now-empty rows 2..6 (oldest-first,
#    matching the existing sheet ordering).
# ---------------------------------------------------------------------------
$newDatesTop = @(44440, 44447, 44454, 44461, 44468)
$newValuesTop = @(277.65, 253.496, 230.328, 308.022, 197.744)

for ($i = 0; $i -lt $newDatesTop.Count; $i++) {
    $row = 2 + $i
    # Copy format from the row that used to be row 2 (now row 7) so the new
    # rows pick up the same date style without introducing new cell styles.
    $ws.Range("A7:B7").Copy($ws.Range("A" + $row + ":B" + $row))
    $ws.Cells.Item($row, 1).Value = $newDatesTop[$i]
    $ws.Cells.Item($row, 2).Value = $newValuesTop[$i]
}

# ---------------------------------------------------------------------------
# 3) Append 2 new observations after the (now shifted) last row, which sits
#    at row 118.
# ---------------------------------------------------------------------------
$newDatesBottom = @(45259, 45266)
$newValuesBottom = @(739.2089999999999, 682.509)

for ($i = 0; $i -lt $newDatesBottom.Count; $i++) {
    $row = 119 + $i
    $ws.Range("A118:B118").Copy($ws.Range("A" + $row + ":B" + $row))
    $ws.Cells.Item($row, 1).Value = $newDatesBottom[$i]
    $ws.Cells.Item($row, 2).Value = $newValuesBottom[$i]
}

# ---------------------------------------------------------------------------
# 4) Refresh the SeriesInfo sheet metadata to match the FRED re-pull.
#    realtime_start / realtime_end / observation_end look like plain dates,
#    so Excel's usual typed-value coercion would turn them into date
#    serials. Write them with a leading apostrophe to force text, then
#    restore the (unformatted) look of the cell via a formats-only paste
#    from a blank scratch cell so the stored style matches the original
#    (no explicit cell style).
# ---------------------------------------------------------------------------
$info = $wb.Worksheets.Item("SeriesInfo")

$info.Range("B3").Value = "'2023-12-08"
$info.Range("B4").Value = "'2023-12-08"
$info.Range("B7").Value = "'2023-12-06"

$info.Range("Z100").Copy()
$info.Range("B3").PasteSpecial(-4122)
$info.Range("B4").PasteSpecial(-4122)
$info.Range("B7").PasteSpecial(-4122)

# last_updated includes a time + UTC offset, so it is not misread as a date.
$info.Range("B14").Value = "2023-12-07 15:34:06-06"

Write-Output "WTREGEN data refreshed"
